$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds text values (e.g. "26.159.78", "0.100") that must
# not be auto-coerced to numbers by Excel, so force Text format first.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.159.78'
$ws.Range("D3").Value = '1.576.42'
$ws.Range("E3").Value = '  -1.71%  '
$ws.Range("E4").Value = '  -0.45%  '
$ws.Range("D5").Value = '208.58'
$ws.Range("E5").Value = '  -1.63%  '
$ws.Range("E6").Value = '  -3.03%  '
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("E8").Value = '  -1.65%  '
$ws.Range("D9").Value = '0.244'
$ws.Range("E9").Value = '  -1.01%  '
$ws.Range("D10").Value = '19.56'
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("D11").Value = '0.0842'
$ws.Range("E11").Value = '  -0.65%  '
$ws.Range("D12").Value = '1.798.12'
$ws.Range("E12").Value = '  -1.71%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.599.88'
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '4.06'
$ws.Range("E14").Value = '  -0.27%  '
$ws.Range("E15").Value = '  -2.18%  '
$ws.Range("D16").Value = '64.39'
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = '26.157.82'
$ws.Range("E17").Value = '  -2.20%  '
$ws.Range("E18").Value = '  -2.02%  '
$ws.Range("D19").Value = '7.26'
$ws.Range("E19").Value = '  +1.09%  '
$ws.Range("D20").Value = '208.66'
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("E22").Value = '  -1.27%  '
$ws.Range("D23").Value = '2.16'
$ws.Range("E23").Value = '  -2.96%  '
$ws.Range("E24").Value = '  -2.30%  '
$ws.Range("D25").Value = '143.81'
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("E27").Value = '  -1.69%  '
$ws.Range("E28").Value = '  -1.92%  '
$ws.Range("E29").Value = '  -0.99%  '
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("E31").Value = '  -1.87%  '
$ws.Range("E32").Value = '  -1.85%  '
$ws.Range("D33").Value = '2.99'
$ws.Range("E33").Value = '  +1.01%  '
$ws.Range("D34").Value = '1.280.32'
$ws.Range("E34").Value = '  -0.64%  '
$ws.Range("E35").Value = '  -1.75%  '
$ws.Range("D36").Value = '0.609'
$ws.Range("E36").Value = '  +3.72%  '
$ws.Range("E37").Value = '  -1.24%  '
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").Value = '1.11'
$ws.Range("E38").Value = '  -9.15%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.0166'
$ws.Range("E39").Value = '  -2.22%  '
$ws.Range("D40").Value = '0.809'
$ws.Range("E40").Value = '  -2.00%  '
$ws.Range("E41").Value = '  -0.42%  '
$ws.Range("E42").Value = '  +2.80%  '
$ws.Range("D43").Value = '0.763'
$ws.Range("E43").Value = '  -1.98%  '
$ws.Range("E44").Value = '  -2.92%  '
$ws.Range("D45").Value = '62.26'
$ws.Range("E45").Value = '  -0.58%  '
$ws.Range("D46").Value = '1.711.25'
$ws.Range("E46").Value = '  -1.69%  '
$ws.Range("D47").Value = '88.67'
$ws.Range("E47").Value = '  -1.93%  '
$ws.Range("E48").Value = '  +0.09%  '
$ws.Range("E49").Value = '  -4.48%  '
$ws.Range("D50").Value = '0.100'
$ws.Range("E50").Value = '  -1.80%  '
$ws.Range("E51").Value = '  -1.52%  '
